$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the used range with a new row 19 by copying the formatting/style
# of the last existing data row (18) down one row. This keeps the new
# row's cell styles consistent with the rest of the data rows (s="2" on
# column A, default style elsewhere) instead of inheriting header formatting
# the way a native Rows.Insert would.
$ws.Range("A18:E18").Copy($ws.Range("A19:E19"))

# Rewrite the full data block (rows 2-19) with the corrected/updated
# forecast series values - this both "shifts" the old rows down by one
# (to make room for the new 2007/2008 observation that used to be missing)
# and refreshes every value per the bugfixed simulation.
$ws.Cells.Item(2,1).Value = 39400
$ws.Cells.Item(2,2).Value = 2007
$ws.Cells.Item(2,3).Value = 0.4235526809466261
$ws.Cells.Item(2,4).Value = 2008
$ws.Cells.Item(2,5).Value = 0.6439341879002525
$ws.Cells.Item(3,1).Value = 39765
$ws.Cells.Item(3,2).Value = 2008
$ws.Cells.Item(3,3).Value = -0.5718076928962645
$ws.Cells.Item(3,4).Value = 2009
$ws.Cells.Item(3,5).Value = -0.1800933741311961
$ws.Cells.Item(4,1).Value = 40130
$ws.Cells.Item(4,2).Value = 2009
$ws.Cells.Item(4,3).Value = 0.3486139762225005
$ws.Cells.Item(4,4).Value = 2010
$ws.Cells.Item(4,5).Value = 0.1555182634501051
$ws.Cells.Item(5,1).Value = 40494
$ws.Cells.Item(5,2).Value = 2010
$ws.Cells.Item(5,3).Value = -0.1384957661262898
$ws.Cells.Item(5,4).Value = 2011
$ws.Cells.Item(5,5).Value = 0.6938817570587785
$ws.Cells.Item(6,1).Value = 40862
$ws.Cells.Item(6,2).Value = 2011
$ws.Cells.Item(6,3).Value = 1.566479473280147
$ws.Cells.Item(6,4).Value = 2012
$ws.Cells.Item(6,5).Value = 0.9614071719361794
$ws.Cells.Item(7,1).Value = 41228
$ws.Cells.Item(7,2).Value = 2012
$ws.Cells.Item(7,3).Value = 0.7307568962936939
$ws.Cells.Item(7,4).Value = 2013
$ws.Cells.Item(7,5).Value = 1.09290550768979
$ws.Cells.Item(8,1).Value = 41592
$ws.Cells.Item(8,2).Value = 2013
$ws.Cells.Item(8,3).Value = 0.818818812164257
$ws.Cells.Item(8,4).Value = 2014
$ws.Cells.Item(8,5).Value = 0.9607602172681418
$ws.Cells.Item(9,1).Value = 41957
$ws.Cells.Item(9,2).Value = 2014
$ws.Cells.Item(9,3).Value = 0.9180054319587239
$ws.Cells.Item(9,4).Value = 2015
$ws.Cells.Item(9,5).Value = 1.375398114243209
$ws.Cells.Item(10,1).Value = 42321
$ws.Cells.Item(10,2).Value = 2015
$ws.Cells.Item(10,3).Value = 1.984684278296656
$ws.Cells.Item(10,4).Value = 2016
$ws.Cells.Item(10,5).Value = 1.473274087935805
$ws.Cells.Item(11,1).Value = 42689
$ws.Cells.Item(11,2).Value = 2016
$ws.Cells.Item(11,3).Value = 1.755995812646982
$ws.Cells.Item(11,4).Value = 2017
$ws.Cells.Item(11,5).Value = 1.681032827388362
$ws.Cells.Item(12,1).Value = 43053
$ws.Cells.Item(12,2).Value = 2017
$ws.Cells.Item(12,3).Value = 1.946965557828384
$ws.Cells.Item(12,4).Value = 2018
$ws.Cells.Item(12,5).Value = 1.755491062323111
$ws.Cells.Item(13,1).Value = 43418
$ws.Cells.Item(13,2).Value = 2018
$ws.Cells.Item(13,3).Value = 1.06432145354225
$ws.Cells.Item(13,4).Value = 2019
$ws.Cells.Item(13,5).Value = 0.776718238020746
$ws.Cells.Item(14,1).Value = 43783
$ws.Cells.Item(14,2).Value = 2019
$ws.Cells.Item(14,3).Value = 1.361817904277696
$ws.Cells.Item(14,4).Value = 2020
$ws.Cells.Item(14,5).Value = 1.316199564471554
$ws.Cells.Item(15,1).Value = 44159
$ws.Cells.Item(15,2).Value = 2020
$ws.Cells.Item(15,3).Value = -4.352425014431304
$ws.Cells.Item(15,4).Value = 2021
$ws.Cells.Item(15,5).Value = 0.03547044462246518
$ws.Cells.Item(16,1).Value = 44525
$ws.Cells.Item(16,2).Value = 2021
$ws.Cells.Item(16,3).Value = -1.761645650979182
$ws.Cells.Item(16,4).Value = 2022
$ws.Cells.Item(16,5).Value = 3.765721202592909
$ws.Cells.Item(17,1).Value = 44890
$ws.Cells.Item(17,2).Value = 2022
$ws.Cells.Item(17,3).Value = 5.20787683103745
$ws.Cells.Item(17,4).Value = 2023
$ws.Cells.Item(17,5).Value = 3.217995704408838
$ws.Cells.Item(18,1).Value = 45254
$ws.Cells.Item(18,2).Value = 2023
$ws.Cells.Item(18,3).Value = -0.9008525709169546
$ws.Cells.Item(18,4).Value = 2024
$ws.Cells.Item(18,5).Value = 0.6027009207580036
$ws.Cells.Item(19,1).Value = 45618
$ws.Cells.Item(19,2).Value = 2024
$ws.Cells.Item(19,3).Value = 0.2738544794132824
$ws.Cells.Item(19,4).Value = 2025
$ws.Cells.Item(19,5).Value = 0.2681899963140832
